$wb = $excel.ActiveWorkbook

# Rename the existing sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ADDCUSTOMER"

# Add a new sheet for OpenAccountTest, placed after sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

# Populate the new sheet data
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "Anne Zimmermann"
$ws2.Range("B2").Value = "Real"

$ws2.Columns.Item(1).ColumnWidth = 16.1666666666667

$ws2.Range("B2").Select()

# Reselect sheet1 and update its selection from E3 -> B6 (sheet1 stays the
# active / tab-selected sheet in the final workbook)
$ws1.Activate()
$ws1.Range("B6").Select()
